$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 58; this shifts rows 58:160 down to 59:161
# and keeps formatting (e.g. the date style on column D) consistent with
# the surrounding rows.
$ws.Rows("58:58").Insert()

# Populate the newly inserted row 58 with the new data record.
$ws.Range("A58").Value = 5
$ws.Range("B58").Value = "Macroferia Regional de Talca"
$ws.Range("C58").Value = "Maule"
$ws.Range("D58").Value = 44935
$ws.Range("E58").Value = 7
$ws.Range("F58").Value = 100112030
$ws.Range("G58").Value = "Poroto granado"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 500
$ws.Range("K58").Value = 38000
$ws.Range("L58").Value = 40000
$ws.Range("M58").Value = 38800
$ws.Range("N58").Value = "$/saco 25 kilos"
$ws.Range("O58").Value = "Región del Maule"
$ws.Range("P58").Value = 1552
$ws.Range("Q58").Value = 25
$ws.Range("R58").Value = "Hortaliza"
